$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.815.56'
$ws.Range('D3').Value = '1.629.28'
$ws.Range('E3').Value = '  -0.63%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.51'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5069'
$ws.Range('E6').Value = '  -0.51%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2579'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06419'
$ws.Range('E9').Value = '  +1.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.33'
$ws.Range('E10').Value = '  -2.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07801'
$ws.Range('E11').Value = '  +0.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.259'
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('D13').Value = '1.629.16'
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('D14').Value = '1.853.89'
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5578'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.25'
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('D17').Value = '0.0₅7517'
$ws.Range('E17').Value = '  -2.91%  '
$ws.Range('D18').Value = '25.820.38'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.78'
$ws.Range('E20').Value = '  -1.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.295'
$ws.Range('E21').Value = '  -2.97%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.794'
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.833'
$ws.Range('E25').Value = '  -3.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1278'
$ws.Range('E26').Value = '  +4.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '140.25'
$ws.Range('E27').Value = '  -1.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.721'
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('E29').Value = '  -1.04%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04854'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.283'
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.184'
$ws.Range('E33').Value = '  -0.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.554'
$ws.Range('E34').Value = '  +0.82%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.380'
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.8939'
$ws.Range('E36').Value = '  -2.23%  '
$ws.Range('D37').Value = '1.133.02'
$ws.Range('E37').Value = '  +3.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.548'
$ws.Range('E38').Value = '  -0.57%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5455'
$ws.Range('E39').Value = '  -1.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01555'
$ws.Range('E40').Value = '  -0.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.0000'
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.574'
$ws.Range('E42').Value = '  -0.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7957'
$ws.Range('E43').Value = '  -1.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '97.16'
$ws.Range('E44').Value = '  -2.04%  '
$ws.Range('D45').Value = '1.777.75'
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('E46').Value = '  -6.81%  '
$ws.Range('E47').Value = '  -2.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.94'
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.619'
$ws.Range('E50').Value = '  +1.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.000'
$ws.Range('E51').Value = '  -0.71%  '
